# SCE_ClientInfo_Template.xlsx - "Added functionality to spinners"
# Update the "# of Students" counts (column C) that are driven by the
# spinner controls, and move the current selection to match the
# author's last cursor position when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Term 1 numbers (rows 2-7: PCOM, BCOM, PM, BA, GL, FS)
$ws.Range("C2").Value = 40
$ws.Range("C3").Value = 30
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 10

# Term 3 numbers (rows 19, 23, 24: BCOM, FS, DXD)
$ws.Range("C19").Value = 20
$ws.Range("C23").Value = 10
$ws.Range("C24").Value = 10

# Move the active selection to reflect where the author left off.
$ws.Range("O21").Select()
